$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-02 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-03 Saturday", 2) | Out-Null
$d.Content.Find.Execute("26×61=1586", $true, $false, $false, $false, $false, $true, 1, $false, "71×12=852", 2) | Out-Null
$d.Content.Find.Execute("25×88=2200", $true, $false, $false, $false, $false, $true, 1, $false, "39×36=1404", 2) | Out-Null
$d.Content.Find.Execute("35×42=1470", $true, $false, $false, $false, $false, $true, 1, $false, "55×15=825", 2) | Out-Null
$d.Content.Find.Execute("32×13=416", $true, $false, $false, $false, $false, $true, 1, $false, "20×53=1060", 2) | Out-Null
$d.Content.Find.Execute("30×64=1920", $true, $false, $false, $false, $false, $true, 1, $false, "14×18=252", 2) | Out-Null
$d.Content.Find.Execute("89×64=5696", $true, $false, $false, $false, $false, $true, 1, $false, "56×60=3360", 2) | Out-Null
$d.Content.Find.Execute("51×60=3060", $true, $false, $false, $false, $false, $true, 1, $false, "79×11=869", 2) | Out-Null
$d.Content.Find.Execute("77×12=924", $true, $false, $false, $false, $false, $true, 1, $false, "56×91=5096", 2) | Out-Null
$d.Content.Find.Execute("38×34=1292", $true, $false, $false, $false, $false, $true, 1, $false, "33×26=858", 2) | Out-Null
$d.Content.Find.Execute("41×22=902", $true, $false, $false, $false, $false, $true, 1, $false, "76×49=3724", 2) | Out-Null
$d.Content.Find.Execute("47×78=3666", $true, $false, $false, $false, $false, $true, 1, $false, "68×50=3400", 2) | Out-Null
$d.Content.Find.Execute("84×49=4116", $true, $false, $false, $false, $false, $true, 1, $false, "30×50=1500", 2) | Out-Null
$d.Content.Find.Execute("37×89=3293", $true, $false, $false, $false, $false, $true, 1, $false, "28×69=1932", 2) | Out-Null
$d.Content.Find.Execute("82×54=4428", $true, $false, $false, $false, $false, $true, 1, $false, "97×54=5238", 2) | Out-Null
$d.Content.Find.Execute("81×72=5832", $true, $false, $false, $false, $false, $true, 1, $false, "14×47=658", 2) | Out-Null
$d.Content.Find.Execute("55×49=2695", $true, $false, $false, $false, $false, $true, 1, $false, "84×15=1260", 2) | Out-Null
$d.Content.Find.Execute("14×21=294", $true, $false, $false, $false, $false, $true, 1, $false, "42×71=2982", 2) | Out-Null
$d.Content.Find.Execute("80×27=2160", $true, $false, $false, $false, $false, $true, 1, $false, "93×23=2139", 2) | Out-Null
$d.Content.Find.Execute("72×51=3672", $true, $false, $false, $false, $false, $true, 1, $false, "46×96=4416", 2) | Out-Null
$d.Content.Find.Execute("27×59=1593", $true, $false, $false, $false, $false, $true, 1, $false, "51×43=2193", 2) | Out-Null
$d.Content.Find.Execute("21×24=504", $true, $false, $false, $false, $false, $true, 1, $false, "16×91=1456", 2) | Out-Null
$d.Content.Find.Execute("28×92=2576", $true, $false, $false, $false, $false, $true, 1, $false, "99×70=6930", 2) | Out-Null
$d.Content.Find.Execute("83×48=3984", $true, $false, $false, $false, $false, $true, 1, $false, "53×57=3021", 2) | Out-Null
$d.Content.Find.Execute("44×51=2244", $true, $false, $false, $false, $false, $true, 1, $false, "90×58=5220", 2) | Out-Null
$d.Content.Find.Execute("15×41=615", $true, $false, $false, $false, $false, $true, 1, $false, "91×20=1820", 2) | Out-Null
